$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Taul1")

# Widen column A to fit the new, longer feature names
$ws.Columns.Item(1).ColumnWidth = 34.6667

# Add the three new backlog rows
$ws.Range("A9").Value = "Ihmisten ryhmittyminen kartalla"
$ws.Range("B9").Value = "Keskitaso"
$ws.Range("C9").Value = "Jos yhdellä alueella on paljon käyttäjiä, käyttäjien kuvakkeet sulautuvat palloon, josta käy ilmi pienellä alueella olevien käyttäjien määrä. "

$ws.Range("A10").Value = "Ryhmään liittyminen ilman käyttäjää"
$ws.Range("B10").Value = "Alhainen"
$ws.Range("C10").Value = "Jos käyttäjä unohtaa kirjautumistietonsa, hän voi liittyä vieraskäyttäjänä ryhmään ryhmän nimen ja salasanan avulla."

$ws.Range("C11").Value = "Layers-valikosta voi piilottaa/näyttää kartan ikonikategorioita (Käyttäjät, tapahtumat, nukkumapaikat…)"
$ws.Range("A11").Value = "Layers-valikko"
$ws.Range("B11").Value = "Korkea"

# Update view state: scroll the window so row 4 is at the top, and
# leave the selection on A18 (mirrors the author's final cursor position)
$ws.Range("A18").Select()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 4
